$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: Prideful Worm ---
$ws.Range("E23").Value = "Slow"

# --- Row 28: Insect Queen (set Extremely High before Very High so the new
#     shared-string indices land in the same order as the target file) ---
$ws.Range("D28").Value = "Extremely High"

# --- Row 24: boss caterpillar ---
$ws.Range("D24").Value = "Very High"

# Apply the same green row formatting (fill + font) used by rows 19-22 to
# rows 23-25 so they match the "customFormat" boss rows above them.
$ws.Range("A23:E23").Interior.Color = $ws.Range("A19:E19").Interior.Color
$ws.Range("A24:E24").Interior.Color = $ws.Range("A19:E19").Interior.Color
$ws.Range("A25:E25").Interior.Color = $ws.Range("A19:E19").Interior.Color

# Row 23 remaining cells
$ws.Range("C23").Value = "Land"
$ws.Range("C23").HorizontalAlignment = $ws.Range("C19").HorizontalAlignment
$ws.Range("C23").VerticalAlignment = $ws.Range("C19").VerticalAlignment
$ws.Range("D23").Value = "High"
$ws.Range("D23").HorizontalAlignment = $ws.Range("D19").HorizontalAlignment

# Row 24 remaining cells
$ws.Range("C24").HorizontalAlignment = $ws.Range("C19").HorizontalAlignment
$ws.Range("C24").VerticalAlignment = $ws.Range("C19").VerticalAlignment
$ws.Range("D24").HorizontalAlignment = $ws.Range("D19").HorizontalAlignment
$ws.Range("E24").Value = "Slow"

# Row 25: Beetle 2
$ws.Range("C25").Value = "Land"
$ws.Range("C25").HorizontalAlignment = $ws.Range("C19").HorizontalAlignment
$ws.Range("C25").VerticalAlignment = $ws.Range("C19").VerticalAlignment
$ws.Range("D25").Value = "Very High"
$ws.Range("D25").HorizontalAlignment = $ws.Range("D19").HorizontalAlignment
$ws.Range("E25").Value = "Normal"

# Row 26: Beetle 3
$ws.Range("D26").Value = "Very High"

# Row 27: Queen Bitter Gourd
$ws.Range("C27").Value = "Land"
$ws.Range("D27").Value = "Very High"
$ws.Range("E27").Value = "Normal"

# Row 28: Insect Queen (remaining cells; D28 already set above)
$ws.Range("C28").Value = "Land"
$ws.Range("E28").Value = "Slow"

# --- View state: rows 23:24 selected, scrolled so row 4 is at the top ---
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
$ws.Rows("23:24").Select()

Write-Host "beetle 2 boss added"
